$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 29.223446
$ws.Range("H2").Value = 87.670338
$ws.Range("I2").Value = 0.0169041244192178
$ws.Range("J2").Value = 0.0169041244192178
$ws.Range("M2").Value = 15.03463666666667
$ws.Range("N2").Value = 45.10391
$ws.Range("O2").Value = 0.2402934356091235
$ws.Range("P2").Value = 0.2402934356091235
$ws.Range("Q2").Value = 439.3638927579533
$ws.Range("R2").Value = 3954.27503482158
$ws.Range("S2").Value = 0.004061950132657925
$ws.Range("T2").Value = 0.004061950132657925

$ws.Range("G3").Value = 29.223446
$ws.Range("H3").Value = 87.670338
$ws.Range("I3").Value = 0.0169041244192178
$ws.Range("J3").Value = 0.0169041244192178
$ws.Range("O3").Value = 0.07715484716443403
$ws.Range("P3").Value = 0.07715484716443403
$ws.Range("Q3").Value = 141.0735749371573
$ws.Range("R3").Value = 1269.662174434416
$ws.Range("S3").Value = 0.001304235136013327
$ws.Range("T3").Value = 0.001304235136013327

$ws.Range("G4").Value = 29.223446
$ws.Range("H4").Value = 87.670338
$ws.Range("I4").Value = 0.0169041244192178
$ws.Range("J4").Value = 0.0169041244192178
$ws.Range("M4").Value = 6.211932333333333
$ws.Range("N4").Value = 18.635797
$ws.Range("O4").Value = 0.09928318157880762
$ws.Range("P4").Value = 0.09928318157880764
$ws.Range("Q4").Value = 181.5340690988207
$ws.Range("R4").Value = 1633.806621889386
$ws.Range("S4").Value = 0.001678295254143957
$ws.Range("T4").Value = 0.001678295254143957

$ws.Range("G5").Value = 29.223446
$ws.Range("H5").Value = 87.670338
$ws.Range("I5").Value = 0.0169041244192178
$ws.Range("J5").Value = 0.0169041244192178
$ws.Range("M5").Value = 36.49384133333334
$ws.Range("N5").Value = 109.481524
$ws.Range("O5").Value = 0.5832685356476348
$ws.Range("P5").Value = 0.5832685356476348
$ws.Range("Q5").Value = 1066.475801537235
$ws.Range("R5").Value = 9598.282213835113
$ws.Range("S5").Value = 0.009859643896402592
$ws.Range("T5").Value = 0.009859643896402592

$ws.Range("I6").Value = 0.9471112884046843
$ws.Range("J6").Value = 0.9471112884046842
$ws.Range("M6").Value = 15.03463666666667
$ws.Range("N6").Value = 45.10391
$ws.Range("O6").Value = 0.2402934356091235
$ws.Range("P6").Value = 0.2402934356091235
$ws.Range("Q6").Value = 24616.86226560192
$ws.Range("R6").Value = 221551.7603904173
$ws.Range("S6").Value = 0.227584625394945
$ws.Range("T6").Value = 0.227584625394945

$ws.Range("I7").Value = 0.9471112884046843
$ws.Range("J7").Value = 0.9471112884046842
$ws.Range("O7").Value = 0.07715484716443403
$ws.Range("P7").Value = 0.07715484716443403
$ws.Range("S7").Value = 0.07307422670457361
$ws.Range("T7").Value = 0.07307422670457361

$ws.Range("I8").Value = 0.9471112884046843
$ws.Range("J8").Value = 0.9471112884046842
$ws.Range("M8").Value = 6.211932333333333
$ws.Range("N8").Value = 18.635797
$ws.Range("O8").Value = 0.09928318157880762
$ws.Range("P8").Value = 0.09928318157880764
$ws.Range("Q8").Value = 10171.06605522043
$ws.Range("R8").Value = 91539.5944969839
$ws.Range("S8").Value = 0.0940322220220207
$ws.Range("T8").Value = 0.0940322220220207

$ws.Range("I9").Value = 0.9471112884046843
$ws.Range("J9").Value = 0.9471112884046842
$ws.Range("M9").Value = 36.49384133333334
$ws.Range("N9").Value = 109.481524
$ws.Range("O9").Value = 0.5832685356476348
$ws.Range("P9").Value = 0.5832685356476348
$ws.Range("Q9").Value = 59752.94817979619
$ws.Range("R9").Value = 537776.5336181658
$ws.Range("S9").Value = 0.5524202142831449
$ws.Range("T9").Value = 0.5524202142831449

$ws.Range("G10").Value = 37.39212666666667
$ws.Range("H10").Value = 112.17638
$ws.Range("I10").Value = 0.02162924801792661
$ws.Range("J10").Value = 0.0216292480179266
$ws.Range("M10").Value = 15.03463666666667
$ws.Range("N10").Value = 45.10391
$ws.Range("O10").Value = 0.2402934356091235
$ws.Range("P10").Value = 0.2402934356091235
$ws.Range("Q10").Value = 562.1770386273112
$ws.Range("R10").Value = 5059.5933476458
$ws.Range("S10").Value = 0.00519736631586941
$ws.Range("T10").Value = 0.00519736631586941

$ws.Range("G11").Value = 37.39212666666667
$ws.Range("H11").Value = 112.17638
$ws.Range("I11").Value = 0.02162924801792661
$ws.Range("J11").Value = 0.0216292480179266
$ws.Range("O11").Value = 0.07715484716443403
$ws.Range("P11").Value = 0.07715484716443403
$ws.Range("Q11").Value = 180.5071511200178
$ws.Range("R11").Value = 1624.56436008016
$ws.Range("S11").Value = 0.001668801325104765
$ws.Range("T11").Value = 0.001668801325104765

$ws.Range("G12").Value = 37.39212666666667
$ws.Range("H12").Value = 112.17638
$ws.Range("I12").Value = 0.02162924801792661
$ws.Range("J12").Value = 0.0216292480179266
$ws.Range("M12").Value = 6.211932333333333
$ws.Range("N12").Value = 18.635797
$ws.Range("O12").Value = 0.09928318157880762
$ws.Range("P12").Value = 0.09928318157880764
$ws.Range("Q12").Value = 232.2773606527622
$ws.Range("R12").Value = 2090.49624587486
$ws.Range("S12").Value = 0.002147420558376872
$ws.Range("T12").Value = 0.002147420558376872

$ws.Range("G13").Value = 37.39212666666667
$ws.Range("H13").Value = 112.17638
$ws.Range("I13").Value = 0.02162924801792661
$ws.Range("J13").Value = 0.0216292480179266
$ws.Range("M13").Value = 36.49384133333334
$ws.Range("N13").Value = 109.481524
$ws.Range("O13").Value = 0.5832685356476348
$ws.Range("P13").Value = 0.5832685356476348
$ws.Range("Q13").Value = 1364.582337689236
$ws.Range("R13").Value = 12281.24103920312
$ws.Range("S13").Value = 0.01261565981857556
$ws.Range("T13").Value = 0.01261565981857556

$ws.Range("G14").Value = 24.817167
$ws.Range("H14").Value = 74.45150100000001
$ws.Range("I14").Value = 0.01435533915817136
$ws.Range("J14").Value = 0.01435533915817136
$ws.Range("M14").Value = 15.03463666666667
$ws.Range("N14").Value = 45.10391
$ws.Range("O14").Value = 0.2402934356091235
$ws.Range("P14").Value = 0.2402934356091235
$ws.Range("Q14").Value = 373.11708894099
$ws.Range("R14").Value = 3358.05380046891
$ws.Range("S14").Value = 0.003449493765651179
$ws.Range("T14").Value = 0.003449493765651179

$ws.Range("G15").Value = 24.817167
$ws.Range("H15").Value = 74.45150100000001
$ws.Range("I15").Value = 0.01435533915817136
$ws.Range("J15").Value = 0.01435533915817136
$ws.Range("O15").Value = 0.07715484716443403
$ws.Range("P15").Value = 0.07715484716443403
$ws.Range("Q15").Value = 119.802656692248
$ws.Range("R15").Value = 1078.223910230232
$ws.Range("S15").Value = 0.001107583998742327
$ws.Range("T15").Value = 0.001107583998742326

$ws.Range("G16").Value = 24.817167
$ws.Range("H16").Value = 74.45150100000001
$ws.Range("I16").Value = 0.01435533915817136
$ws.Range("J16").Value = 0.01435533915817136
$ws.Range("M16").Value = 6.211932333333333
$ws.Range("N16").Value = 18.635797
$ws.Range("O16").Value = 0.09928318157880762
$ws.Range("P16").Value = 0.09928318157880764
$ws.Range("Q16").Value = 154.162562109033
$ws.Range("R16").Value = 1387.463058981297
$ws.Range("S16").Value = 0.001425243744266095
$ws.Range("T16").Value = 0.001425243744266095

$ws.Range("G17").Value = 24.817167
$ws.Range("H17").Value = 74.45150100000001
$ws.Range("I17").Value = 0.01435533915817136
$ws.Range("J17").Value = 0.01435533915817136
$ws.Range("M17").Value = 36.49384133333334
$ws.Range("N17").Value = 109.481524
$ws.Range("O17").Value = 0.5832685356476348
$ws.Range("P17").Value = 0.5832685356476348
$ws.Range("Q17").Value = 905.6737548408361
$ws.Range("R17").Value = 8151.063793567525
$ws.Range("S17").Value = 0.00837301764951176
$ws.Range("T17").Value = 0.00837301764951176
